# Add columns I (I0) and J (IF) to the sheet, populating header and data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy H1's formatting (bold/border/alignment) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I column is always 1, J column mirrors column H (IP)
for ($r = 2; $r -le 10; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
